$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.634
$ws.Range("C3").Value = -12.634
$ws.Range("E4").Value = 12.662
$ws.Range("C5").Value = -12.634
$ws.Range("E6").Value = 12.336
$ws.Range("D7").Value = -7.237
$ws.Range("A9").Value = -20.775
$ws.Range("D9").Value = -7.873
$ws.Range("E10").Value = 12.436
$ws.Range("C11").Value = -12.628
$ws.Range("E11").Value = 13.072
$ws.Range("C12").Value = -12.628
$ws.Range("A13").Value = -21.99
$ws.Range("A16").Value = -20.911
$ws.Range("A18").Value = -21.577
$ws.Range("A20").Value = -21.753
$ws.Range("C21").Value = -11.964
$ws.Range("D21").Value = -7.8
$ws.Range("E21").Value = 13.017
$ws.Range("E25").Value = 12.599
